$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.829.56'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.874.39'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.75%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.08%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5365'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3746'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07196'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.52%  '
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8900'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08180'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.875.53'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '93.36'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.99%  '
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.84'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008529'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.91%  '
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.865.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.991'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.62'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.403'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.292'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.07'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.718'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.623'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09134'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8103'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05014'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.175'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.948'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6045'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.226'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.634'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.49%  '
$ws.Range('E39').Value = '  -2.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.071'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.625'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('B43').Value = 'Decentraland'
$ws.Range('C43').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5138'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.85%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '115.14'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1499'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.65%  '
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.646'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.993'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06082'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '62.24'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.01%  '
